# "1st iteration (#141)" - regenerated build metadata sheet:
#  - Metadata!B7 ("Experimental") gets a value: text "true"
#  - Metadata!B8 ("Date") is refreshed to the new build timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Leading apostrophe forces literal text (not the Boolean TRUE) -
# matches the shared-string "true" added by the diff.
$ws.Range("B7").Value = "'true"

$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
